$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell (H1) onto the new
# header cells so they pick up the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-16 for columns I (I0) and J (IF)
$data = @(
    @(6, 7),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(7, 9),
    @(9, 9),
    @(7, 8),
    @(8, 9),
    @(4, 5),
    @(9, 9)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
